# Automatische test-sync: 2025-06-27 22:38:50
# Adds the newest "Wanneer zijn jullie open?" test-mail log entry as row 9
# on the "Logs" sheet, extends the conditional formatting ranges to cover
# it, and bumps the matching category counter on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$antwoord = "Beste klant,`nBedankt voor uw e-mail. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. We zijn gesloten in het weekend. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf]"

$ws.Cells.Item(9, 1).Value = "Wanneer zijn jullie open?"
$ws.Cells.Item(9, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item(9, 3).Value = "Testmail #1: Wanneer zijn jullie open?"
$ws.Cells.Item(9, 4).Value = "Openingstijden / Locatie"
$ws.Cells.Item(9, 5).Value = $antwoord
$ws.Cells.Item(9, 6).Value = "2025-06-27 22:38:47"
$ws.Cells.Item(9, 7).Value = "Ja"
$ws.Cells.Item(9, 8).Value = "Nee"
$ws.Cells.Item(9, 9).Value = "Ja"

# The new row keeps the sheet's default row height (no custom height).
$ws.Rows.Item(9).RowHeight = $ws.Rows.Item(8).RowHeight

# Extend the conditional-formatting ranges so they keep covering column
# D/G/H/I down to the newly-added row 9.
$ws.Range("D2:D8").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D9"))
$ws.Range("G2:G8").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G9"))
$ws.Range("H2:H8").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H9"))
$ws.Range("I2:I8").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I9"))

# Bump the "Openingstijden / Locatie" tally on the Dashboard sheet.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(2, 2).Value = 6
